# New benchmarks after caching
# Append two new benchmark rows (35 & 36) to the "Development" sheet,
# recording the run after introducing caching, and move the active
# selection onto the newly entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Development")

# Seed row 35/36 from row 34 first so the date cells inherit the same
# date style (s="1") as every other row in column A, instead of Excel
# minting a brand new cell style.
[void]$ws.Range("A34").Copy($ws.Range("A35"))
[void]$ws.Range("A34").Copy($ws.Range("A36"))

# These two rows are a little taller than the rest of the table.
$ws.Rows.Item(35).RowHeight = 15
$ws.Rows.Item(36).RowHeight = 15

# Row 35 - first run after introducing caching
$ws.Range("A35").Value = 44539
$ws.Range("B35").Value = "7.0.1 (develop)"
$ws.Range("C35").Value = 364.14
$ws.Range("D35").Value = 7.42
$ws.Range("E35").Value = 307
$ws.Range("F35").Value = 120
$ws.Range("G35").Value = "Linux"
$ws.Range("H35").Value = "Intel skylake (core m7) Intel Core(TM)i7-6500U@2.50GHz"
$ws.Range("L35").Value = "Default settings, after introducing caching"

# Row 36 - second run after introducing caching
$ws.Range("A36").Value = 44539
$ws.Range("B36").Value = "7.0.1 (develop)"
$ws.Range("C36").Value = 332
$ws.Range("D36").Value = 7.3
$ws.Range("E36").Value = 347
$ws.Range("F36").Value = 97.4
$ws.Range("G36").Value = "Linux"
$ws.Range("H36").Value = "Intel skylake (core m7) Intel Core(TM)i7-6500U@2.50GHz"
$ws.Range("L36").Value = "Default settings, after introducing caching"

# Leave the cursor where the author left it after typing the new data.
[void]$ws.Range("C36").Select()
